$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, pushing existing rows 133..307 down to 134..308.
$ws.Rows(133).Insert()

# Populate the newly inserted row 133 with a fresh weekly observation
# (same market/category metadata as every other row in this sheet).
$ws.Range("A133").Value = 3
$ws.Range("B133").Value = "Femacal de La Calera"
$ws.Range("C133").Value = "Coquimbo"
$ws.Range("D133").Value = 44679
$ws.Range("E133").Value = 5
$ws.Range("F133").Value = 100112039
$ws.Range("G133").Value = "Ciboulette"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 180
$ws.Range("K133").Value = 1500
$ws.Range("L133").Value = 1500
$ws.Range("M133").Value = 1500
$ws.Range("N133").Value = '$/docena de atados'
$ws.Range("O133").Value = "Provincia de Quillota"
$ws.Range("P133").Value = 500
$ws.Range("Q133").Value = 3
$ws.Range("R133").Value = "Hortaliza"
